$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.89%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.44%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.879"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.38%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'0.00%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.934"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.03%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.285"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'34.05%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8744"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.64%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'5.39%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05055"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.63%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'2.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02956"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.65%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09068"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001582"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.82%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006323"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.54%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006038"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'4.59%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'3.323"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.50%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.284"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.49%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.30%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1335"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.73%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.932"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.97%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04372"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.18%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001174"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.75%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001617"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-4.37%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.04095"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.16%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007018"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.35%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.86%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-0.45%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-13.10%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'1.58%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-11.17%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.486"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-37.51%"
$ws.Range("E47").Style = "Normal"
